# edit.ps1 - applies the "Add files via upload" edit to TASK-2.docx
#
# Summary of content changes (see diff):
#   1. Insert a new paragraph "1." right before "CREATE TABLE MATCHES("
#   2. Insert a new paragraph "2." right before "CREATE TABLE DELIVERIES("
#   3. After the DELIVERIES table's closing ");", turn the two blank
#      paragraphs / add new ones so we get:
#         3.  FILE UPLOAD USING INBUILT UTILITIES.
#         (blank)
#         4.  FILE UPLOAD USING INBUILT UTILITIES.
#         (blank)
#   4. Prefix the trailing SELECT statements with 5. through 10.

$d = $word.ActiveDocument

function Get-ParaIndexForText($searchText) {
    $rng = $d.Content
    $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $rng.Paragraphs(1).Index
}

# ---------------------------------------------------------------------
# 1. "1." before "CREATE TABLE MATCHES("
# ---------------------------------------------------------------------
$idx = Get-ParaIndexForText("CREATE TABLE MATCHES(")
$d.Paragraphs($idx).Range.InsertParagraphBefore()
$d.Paragraphs($idx).Range.Text = "1."

# ---------------------------------------------------------------------
# 2. "2." before "CREATE TABLE DELIVERIES("
# ---------------------------------------------------------------------
$idx = Get-ParaIndexForText("CREATE TABLE DELIVERIES(")
$d.Paragraphs($idx).Range.InsertParagraphBefore()
$d.Paragraphs($idx).Range.Text = "2."

# ---------------------------------------------------------------------
# 3. Restructure the two blank paragraphs before
#    "SELECT * FROM IPL_DELIVERIES;" into the "3./4." numbered items.
# ---------------------------------------------------------------------
$idxSelDeliv = Get-ParaIndexForText("SELECT * FROM IPL_DELIVERIES;")
$idxBlank1 = $idxSelDeliv - 2
$d.Paragraphs($idxBlank1).Range.Text = "3.  FILE UPLOAD USING INBUILT UTILITIES."

# Insert paragraph "4. ..." plus a trailing blank paragraph right before
# the "SELECT * FROM IPL_DELIVERIES;" paragraph (after the 2nd, still
# blank, paragraph).
$idxSelDeliv = Get-ParaIndexForText("SELECT * FROM IPL_DELIVERIES;")
$d.Paragraphs($idxSelDeliv).Range.InsertParagraphBefore()
$idxSelDeliv = Get-ParaIndexForText("SELECT * FROM IPL_DELIVERIES;")
$d.Paragraphs($idxSelDeliv - 1).Range.Text = "4.  FILE UPLOAD USING INBUILT UTILITIES."

$idxSelDeliv = Get-ParaIndexForText("SELECT * FROM IPL_DELIVERIES;")
$d.Paragraphs($idxSelDeliv).Range.InsertParagraphBefore()

# ---------------------------------------------------------------------
# 4. Number the remaining SELECT statements (5. through 10.)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("SELECT * FROM IPL_DELIVERIES;", $false, $false, $false, $false, $false, $true, 1, $false, "5. SELECT * FROM IPL_DELIVERIES;", 2) | Out-Null
$d.Content.Find.Execute("SELECT * FROM IPL_MATCHES;", $false, $false, $false, $false, $false, $true, 1, $false, "6. SELECT * FROM IPL_MATCHES;", 2) | Out-Null
$d.Content.Find.Execute("SELECT * FROM IPL_MATCHES WHERE DATEOFMATCH='02-05-2013';", $false, $false, $false, $false, $false, $true, 1, $false, "7. SELECT * FROM IPL_MATCHES WHERE DATEOFMATCH='02-05-2013';", 2) | Out-Null
$d.Content.Find.Execute("SELECT * FROM IPL_MATCHES WHERE RESULT_MARGIN>100;", $false, $false, $false, $false, $false, $true, 1, $false, "8. SELECT * FROM IPL_MATCHES WHERE RESULT_MARGIN>100;", 2) | Out-Null
$d.Content.Find.Execute("SELECT * FROM IPL_MATCHES WHERE RESULT='tie' ORDER BY(ID);", $false, $false, $false, $false, $false, $true, 1, $false, "9. SELECT * FROM IPL_MATCHES WHERE RESULT='tie' ORDER BY(ID);", 2) | Out-Null
$d.Content.Find.Execute("SELECT CITY,COUNT(CITY) FROM IPL_MATCHES GROUP BY(CITY);", $false, $false, $false, $false, $false, $true, 1, $false, "10. SELECT CITY,COUNT(CITY) FROM IPL_MATCHES GROUP BY(CITY);", 2) | Out-Null

Write-Output "done"
